$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "C5.101"
$ws.Range("A3").Value = "C5.102"

$rooms = @("C5.103","C5.104","C5.105","C5.106","C5.107","C5.108","C5.109","C5.110","C5.111","C5.112")
$row = 4
foreach ($room in $rooms) {
    $ws.Cells.Item($row, 1).Value = $room
    $row = $row + 1
}

$ws.Range("A14").Select()
